$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.599.34"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").Value = "1.593.79"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("E6").Value = "  +0.80%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -0.98%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.40"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.08%  "
$ws.Range("E11").Value = "  +0.10%  "
$ws.Range("D12").Value = "1.818.16"
$ws.Range("E12").Value = "  +0.34%  "
$ws.Range("D13").Value = "1.609.74"
$ws.Range("E13").Value = "  +1.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.56%  "
$ws.Range("D17").Value = "26.584.21"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").Value = "0.0₃0729"
$ws.Range("E18").Value = "  +0.34%  "
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "207.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.54%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.30"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.14"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.08%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.71%  "
$ws.Range("E28").Value = "  +0.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.20"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.64%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0503"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.61%  "
$ws.Range("E31").Value = "  +0.14%  "
$ws.Range("E32").Value = "  -0.20%  "
$ws.Range("E33").Value = "  -1.59%  "
$ws.Range("E34").Value = "  +0.51%  "
$ws.Range("D35").Value = "1.280.18"
$ws.Range("E35").Value = "  -1.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.45"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.86%  "
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("E38").Value = "  -0.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.838"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.06%  "
$ws.Range("E40").Value = "  +0.15%  "
$ws.Range("E41").Value = "  +1.44%  "
$ws.Range("E42").Value = "  +1.36%  "
$ws.Range("E43").Value = "  -0.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.77"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.917"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.51%  "
$ws.Range("D46").Value = "1.729.78"
$ws.Range("E46").Value = "  +0.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "89.38"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.43%  "
$ws.Range("E48").Value = "  -1.37%  "
$ws.Range("D49").Value = "0.0₆0104"
$ws.Range("E49").Value = "  -1.15%  "
$ws.Range("E50").Value = "  +3.82%  "
$ws.Range("E51").Value = "  +0.65%  "
